$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FrameCounts")

# Fix existing rows 53/54 (B column corrections)
$ws.Range("B53").Value = 7144
$ws.Range("B54").Value = 7174

# Fill in newly tracked rows 55-63
$ws.Range("A55").Value = "X = 282"
$ws.Range("B55").Value = 7331
$ws.Range("C55").Value = 7647

$ws.Range("A56").Value = "Batman appears screen 2"
$ws.Range("B56").Value = 8507
$ws.Range("C56").Value = 8823

$ws.Range("A57").Value = "X = 176"
$ws.Range("B57").Value = 8588
$ws.Range("C57").Value = 8912

$ws.Range("A58").Value = "X = 299"
$ws.Range("B58").Value = 8648
$ws.Range("C58").Value = 8974

$ws.Range("A59").Value = "Begin walljump"
$ws.Range("B59").Value = 9079
$ws.Range("C59").Value = 9639

$ws.Range("A60").Value = "Black screen"
$ws.Range("B60").Value = 9164
$ws.Range("C60").Value = 9724

$ws.Range("A61").Value = "HP = 26"
$ws.Range("B61").Value = 9320
$ws.Range("C61").Value = 9857

$ws.Range("A62").Value = "HP = 0"
$ws.Range("B62").Value = 9468
$ws.Range("C62").Value = 9997

$ws.Range("A63").Value = "Batman disappears"
$ws.Range("B63").Value = 9810
$ws.Range("C63").Value = 10340

# Apply the same style (borders) used by the rest of column A/B/C in this table
$ws.Range("A55:C63").Style = $ws.Range("A54").Style

# Update the frozen pane / view position and active selection to match the new scroll position
$ws.Application.ActiveWindow.ScrollRow = 51
$ws.Range("C64").Select()
